$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New gene-position values produced after switching the upstream filter to
# the rRNA feature type. These fill A2:A194 (193 values), replacing the
# previously-empty cells and extending the sheet with 77 new rows.
$values = @(
    "1065828",
    "1533026",
    "1065657",
    "133545",
    "1532999",
    "1537165",
    "137512",
    "116420",
    "1068403",
    "137410",
    "1537156",
    "1069669",
    "1243701",
    "742246",
    "1535549",
    "1066422",
    "112798",
    "1071450",
    "1080147",
    "113832",
    "113820",
    "1065882",
    "1076064",
    "436520",
    "110054",
    "288650",
    "1066155",
    "1536054",
    "1532695",
    "116738",
    "1536070",
    "110187",
    "398086",
    "133047",
    "137700",
    "133795",
    "137704",
    "1533001",
    "1065803",
    "282889",
    "1535989",
    "1071675",
    "356961",
    "132815",
    "418736",
    "116677",
    "137674",
    "132828",
    "1061972",
    "1061247",
    "1071801",
    "115705",
    "1071797",
    "1071726",
    "1065973",
    "1076417",
    "1066125",
    "354965",
    "1069550",
    "138415",
    "1536649",
    "1532844",
    "132833",
    "1536067",
    "1066376",
    "1068384",
    "1535545",
    "1068377",
    "1536069",
    "137565",
    "354884",
    "1069575",
    "135135",
    "1066130",
    "133523",
    "1533007",
    "137403",
    "1529742",
    "133524",
    "109788",
    "1071725",
    "137587",
    "1533076",
    "436521",
    "742239",
    "1061325",
    "1536445",
    "354923",
    "421434",
    "133891",
    "132432",
    "109787",
    "1061727",
    "115729",
    "138418",
    "138541",
    "1065729",
    "1065783",
    "1532990",
    "1065699",
    "135697",
    "116640",
    "1533061",
    "133640",
    "1529754",
    "1071645",
    "1529738",
    "1354782",
    "1535983",
    "1065885",
    "135688",
    "116406",
    "1071795",
    "1535548",
    "1061682",
    "1060477",
    "132858",
    "133693",
    "354948",
    "1541286",
    "1532996",
    "1536145",
    "1535541",
    "132827",
    "1068399",
    "1535542",
    "132829",
    "1061457",
    "137314",
    "1536424",
    "1536141",
    "133802",
    "1537075",
    "356954",
    "1536214",
    "132784",
    "132823",
    "1425591",
    "1536247",
    "1065996",
    "288697",
    "1066064",
    "1061919",
    "109824",
    "109661",
    "132866",
    "1065861",
    "1071639",
    "1061346",
    "1536886",
    "132864",
    "109658",
    "1536139",
    "132436",
    "1065782",
    "1536730",
    "367304",
    "1536246",
    "1066119",
    "133547",
    "1060471",
    "355004",
    "356962",
    "1060623",
    "1066377",
    "357283",
    "1063847",
    "742245",
    "1536004",
    "1061973",
    "1535539",
    "1071723",
    "1071660",
    "1062014",
    "1537140",
    "424500",
    "1065726",
    "138424",
    "1533074",
    "1068357",
    "436338",
    "137277",
    "1536019",
    "1071387",
    "138578",
    "1066043",
    "1061628",
    "1066128",
    "109657",
    "356957",
    "1664188",
    "116644",
    "1036498"
)

$n = $values.Count
$arr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $arr[$i,0] = $values[$i]
}

# A2:A194 -> all 193 new values in one shot (A2 is row index 2, i.e. right
# after the header row A1).
$ws.Range("A2:A194").Value = $arr

# Rows 118-194 (the extra rRNA rows beyond the previous 117-row extent) pick
# up the bold "group header" look used by A1, instead of the plain style
# that rows 2-117 already have.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A118:A194").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the saved selection/scroll state from the edited workbook.
$ws.Range("A2:A194").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 171
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "Filled A2:A194 with $n values and restyled A118:A194"
